$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")
$ws.Activate()

# Fix seeding error: the assumed number of coital acts per week changed
# from 5 to 4.5 (peak coital acts/partnership = 52 weeks * acts/week).
$ws.Range("AI2").Formula = "=52*4.5"
$ws.Range("AI3").Formula = "=52*4.5"

# Update the active selection/scroll position left behind by the edit.
$ws.Range("AI3").Select()
